# Updates crypto price/volume data per the Dec 27 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) are stored as plain text in the sheet (prices
# use "." as a thousands separator in some rows, e.g. "42.853.03", and volumes
# keep padding spaces, e.g. "  -0.24%  "). Force text format before writing so
# Excel does not reinterpret numeric-looking strings as numbers, then clear the
# temporary formatting so the cell style is left exactly as it was originally.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.853.03"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.24%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.252.07"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "113.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "295.17"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.86%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -0.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.21%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.32"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0925"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.44"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.96"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.07"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +22.16%  "

$ws.Range("E15").Value = "  -1.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.12"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.587.80"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.245.85"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.814.22"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.35%  "

$ws.Range("E20").Value = "  +6.92%  "

$ws.Range("E21").Value = "  -0.81%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.52"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.37"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +11.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.98%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "250.51"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +8.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.36%  "

$ws.Range("E27").Value = "  -0.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.57"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.48%  "

$ws.Range("E29").Value = "  -1.21%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.76"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.83%  "

$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.98"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.97"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.13"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0889"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.99%  "

$ws.Range("E35").Value = "  +2.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.09"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +9.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.28"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.128"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.06%  "

$ws.Range("E39").Value = "  +0.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.105"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.44"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.46"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.46%  "

$ws.Range("E43").Value = "  -0.75%  "

$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.54"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.33"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.51"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.49%  "

$ws.Range("E48").Value = "  +2.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.85"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.63"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.46%  "

$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.87"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.16%  "
